$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed "Price" figures in column D are plain text in the source feed
# (e.g. trailing zeros such as "0.3790", or thousand-grouped values such as
# "28.105.90" that are not valid numbers anyway). Force the cells that would
# otherwise be auto-recognised as numbers into Text format before typing the
# new values so the stored text is preserved exactly, then clear the format
# again so the cells end up with the same (default/general) style as before -
# only their content changes.
$textRanges = @(
    "D5",
    "D7:D11",
    "D13:D15",
    "D17",
    "D19:D20",
    "D22",
    "D24",
    "D26:D27",
    "D29:D34",
    "D36",
    "D38:D42",
    "D44:D51"
)
foreach ($a in $textRanges) {
    $ws.Range($a).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.105.90"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "1.800.02"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "316.75"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.5438"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").Value = "0.3790"
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").Value = "0.07470"
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").Value = "41.95"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "1.093"
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "6.215"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "20.42"
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("D15").Value = "7.371"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "1.796.66"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "89.30"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("D19").Value = "0.06520"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "17.40"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "5.937"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "28.136.81"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("D24").Value = "11.18"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Value = "156.18"
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").Value = "20.38"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "2.006.76"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "2.336"
$ws.Range("E29").Value = "  -1.81%  "
$ws.Range("D30").Value = "121.79"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").Value = "0.1104"
$ws.Range("E31").Value = "  +7.12%  "
$ws.Range("D32").Value = "1.113"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").Value = "3.672"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "5.539"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("E35").Value = "  +5.90%  "
$ws.Range("D36").Value = "0.2212"
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("D38").Value = "5.066"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "8.435"
$ws.Range("E39").Value = "  -5.14%  "
$ws.Range("D40").Value = "11.16"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").Value = "0.6158"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").Value = "1.174"
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("D44").Value = "13.39"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "3.686"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "0.5733"
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("D47").Value = "124.48"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").Value = "1.182"
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("D49").Value = "1.916"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").Value = "0.06816"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("D51").Value = "0.00000000300"
$ws.Range("E51").Value = "  +40.02%  "

foreach ($a in $textRanges) {
    $ws.Range($a).ClearFormats()
}
